$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.266.25"
$ws.Range("E2").Value = "  -0.94%  "
$ws.Range("D3").Value = "3.351.42"
$ws.Range("E3").Value = "  -1.97%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "571.53"
$ws.Range("E5").Value = "  -1.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.26"
$ws.Range("E6").Value = "  -2.40%  "
$ws.Range("D8").Value = "3.351.03"
$ws.Range("E8").Value = "  -1.95%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.469"
$ws.Range("E9").Value = "  -1.68%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.46"
$ws.Range("E10").Value = "  -3.34%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.123"
$ws.Range("E11").Value = "  -2.62%  "
$ws.Range("E12").Value = "  -1.88%  "
$ws.Range("D13").Value = "3.928.88"
$ws.Range("E13").Value = "  -1.66%  "
$ws.Range("E14").Value = "  +1.66%  "
$ws.Range("E15").Value = "  -3.70%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.83"
$ws.Range("E16").Value = "  +0.56%  "
$ws.Range("D17").Value = "3.359.23"
$ws.Range("E17").Value = "  -1.47%  "
$ws.Range("D18").Value = "61.384.98"
$ws.Range("E18").Value = "  -0.78%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.92"
$ws.Range("E19").Value = "  -1.45%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.85"
$ws.Range("E20").Value = "  -1.90%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.28"
$ws.Range("E21").Value = "  -2.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "374.79"
$ws.Range("E22").Value = "  -4.35%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.550"
$ws.Range("E23").Value = "  -4.33%  "
$ws.Range("D24").Value = "3.507.07"
$ws.Range("E24").Value = "  -1.40%  "
$ws.Range("E25").Value = "  -0.14%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "70.88"
$ws.Range("E26").Value = "  -0.54%  "
$ws.Range("E27").Value = "  -2.65%  "
$ws.Range("E28").Value = "  +6.84%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.47"
$ws.Range("E29").Value = "  -4.63%  "
$ws.Range("E30").Value = "  +0.42%  "
$ws.Range("E31").Value = "  +3.21%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.14"
$ws.Range("E32").Value = "  -2.41%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.14"
$ws.Range("E33").Value = "  -0.81%  "
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.50"
$ws.Range("E35").Value = "  -0.93%  "
$ws.Range("E36").Value = "  -6.39%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.78"
$ws.Range("E37").Value = "  -4.15%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.52"
$ws.Range("E38").Value = "  -2.52%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "164.63"
$ws.Range("E39").Value = "  +1.53%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0764"
$ws.Range("E40").Value = "  -4.68%  "
$ws.Range("E41").Value = "  +0.18%  "
$ws.Range("E42").Value = "  -1.36%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.769"
$ws.Range("E43").Value = "  -1.21%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.20"
$ws.Range("E44").Value = "  -2.20%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.30"
$ws.Range("E45").Value = "  +0.09%  "
$ws.Range("E46").Value = "  -3.28%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "24.05"
$ws.Range("E47").Value = "  +2.08%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.83"
$ws.Range("E48").Value = "  -2.90%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "22.79"
$ws.Range("E49").Value = "  -1.17%  "
$ws.Range("D50").Value = "2.361.18"
$ws.Range("E50").Value = "  -0.07%  "
$ws.Range("E51").Value = "  -2.51%  "
